# orcerc_model.xlsx update
# - Processes sheet ("Hoja4"): move the "type" column (old column D) to be right
#   after "key" (new column B), shifting fuel/product right; resize columns;
#   change active selection to E1:E8 and make this the active tab.
# - WasteDefinition sheet ("Hoja9"): no longer the active tab (handled
#   automatically by activating Processes instead).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processes")

# Move column D ("type") so that it becomes column B; fuel/product shift
# right into C/D. This is the same as selecting column D, Cut, and doing
# "Insert Cut Cells" in front of column B.
$ws.Columns.Item(4).Cut()
$ws.Columns.Item(2).Insert()

# Resize the columns to their final widths (Excel ColumnWidth is in
# characters; stored sheet width = ColumnWidth + 5/6).
$ws.Columns.Item(2).ColumnWidth = 12.307291666666666   # B - type
$ws.Columns.Item(3).ColumnWidth = 10.022135416666666   # C - fuel
$ws.Columns.Item(4).ColumnWidth = 9.592447916666666    # D - product
$ws.Columns.Item(5).ColumnWidth = 16.022135416666668   # E - description
$ws.Columns.Item(6).ColumnWidth = 15.877604166666666   # F
$ws.Columns.Item(7).ColumnWidth = 9.166666666666666    # G

# Make Processes the active sheet/tab and update its selection.
$ws.Activate()
$ws.Range("E1:E8").Select()
